# Update 2023 06 (june) 06
# - add a new row (C020_FT Dataprep_Base.ipynb) to the "code" sheet
# - drop the (accidental) bold weight that had crept onto the E4:G6 wrap cells
# - move the active selection to G9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Tidy up existing formatting -----------------------------------
# E4:G6 currently carry a bold font together with the wrap/border/centre
# alignment; the sibling rows (E7:H8) use the same wrap/border/centre
# alignment WITHOUT bold. Copy that non-bold wrapped format across so the
# unused "bold + wrap" style collapses away.
$ws.Range("E7").Copy()
$ws.Range("E4:G6").PasteSpecial(-4122)   # xlPasteFormats

# C9 currently only has a border (no centring); give it the same
# border+centre look used everywhere else in the table (copy from D9).
$ws.Range("D9").Copy()
$ws.Range("C9").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = 0

# --- 2. Fill in the new table row --------------------------------------
$ws.Range("C9").Value = "C020_FT Dataprep_Base.ipynb`r`n"
$ws.Range("D9").Value = "/code-cloud/"
$ws.Range("E9").Value = "raw_base_2023-06-05.csv"

# G9 needs the wrap/border/centre style (same family as E7/E8) plus the
# description text.
$ws.Range("E7").Copy()
$ws.Range("G9").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G9").Value = "Takes the raw data and adds :`r`na) winner, top3, measurement for random `r`nb) winner, top3, measurement for fast-track cols `r`nthen filters for the created columns and exports. `r`nMeasurement variable include (hits, profitability for winner, all place to wins, lay last)`r`n"

# Row grew tall because of the wrapped description - match the authored height.
$ws.Rows.Item(9).RowHeight = 144

# --- 3. Selection moved to G9 in the saved file ------------------------
$ws.Range("G9").Select() | Out-Null

Write-Output "applied"
